$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "agnihotriaman@gmail.com"
$ws.Range("B2").Value = "124ef1"
$ws.Range("C2").Value = "{1: 'e_que_1', 2: 'm_stck_8', 3: 'm_ll_7', 4: 'm_que_7', 5: 'm_str_2', 6: 'e_que_3', 7: 'm_ll_9', 8: 'e_que_7'}"
